{"js": "const pairs = [\n  [\"2024-06-06 Thursday\", \"2024-06-07 Friday\"],\n  [\"76\u00d743=3268\", \"18\u00d711=198\"],\n  [\"26\u00d793=2418\", \"38\u00d774=2812\"],\n  [\"37\u00d752=1924\", \"24\u00d781=1944\"],\n  [\"59\u00d747=2773\", \"18\u00d763=1134\"],\n  [\"21\u00d798=2058\", \"72\u00d717=1224\"],\n  [\"49\u00d755=2695\", \"16\u00d727=432\"],\n  [\"90\u00d744=3960\", \"81\u00d778=6318\"],\n  [\"39\u00d771=2769\", \"56\u00d788=4928\"],\n  [\"50\u00d746=2300\", \"31\u00d739=1209\"],\n  [\"97\u00d724=2328\", \"20\u00d712=240\"],\n  [\"71\u00d764=4544\", \"98\u00d787=8526\"],\n  [\"61\u00d791=5551\", \"22\u00d769=1518\"],\n  [\"92\u00d740=3680\", \"80\u00d790=7200\"],\n  [\"62\u00d775=4650\", \"53\u00d721=1113\"],\n  [\"87\u00d777=6699\", \"13\u00d773=949\"],\n  [\"86\u00d716=1376\", \"67\u00d742=2814\"],\n  [\"77\u00d736=2772\", \"58\u00d793=5394\"],\n  [\"13\u00d787=1131\", \"88\u00d731=2728\"],\n  [\"87\u00d714=1218\", \"51\u00d733=1683\"],\n  [\"79\u00d759=4661\", \"73\u00d726=1898\"],\n  [\"46\u00d762=2852\", \"77\u00d796=7392\"],\n  [\"29\u00d743=1247\", \"95\u00d759=5605\"],\n  [\"16\u00d775=1200\", \"77\u00d730=2310\"],\n  [\"24\u00d722=528\", \"52\u00d714=728\"],\n  [\"19\u00d771=1349\", \"65\u00d780=5200\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2024-06-06 Thursday\", \"2024-06-07 Friday\")\n    ,@(\"76\u00d743=3268\", \"18\u00d711=198\")\n    ,@(\"26\u00d793=2418\", \"38\u00d774=2812\")\n    ,@(\"37\u00d752=1924\", \"24\u00d781=1944\")\n    ,@(\"59\u00d747=2773\", \"18\u00d763=1134\")\n    ,@(\"21\u00d798=2058\", \"72\u00d717=1224\")\n    ,@(\"49\u00d755=2695\", \"16\u00d727=432\")\n    ,@(\"90\u00d744=3960\", \"81\u00d778=6318\")\n    ,@(\"39\u00d771=2769\", \"56\u00d788=4928\")\n    ,@(\"50\u00d746=2300\", \"31\u00d739=1209\")\n    ,@(\"97\u00d724=2328\", \"20\u00d712=240\")\n    ,@(\"71\u00d764=4544\", \"98\u00d787=8526\")\n    ,@(\"61\u00d791=5551\", \"22\u00d769=1518\")\n    ,@(\"92\u00d740=3680\", \"80\u00d790=7200\")\n    ,@(\"62\u00d775=4650\", \"53\u00d721=1113\")\n    ,@(\"87\u00d777=6699\", \"13\u00d773=949\")\n    ,@(\"86\u00d716=1376\", \"67\u00d742=2814\")\n    ,@(\"77\u00d736=2772\", \"58\u00d793=5394\")\n    ,@(\"13\u00d787=1131\", \"88\u00d731=2728\")\n    ,@(\"87\u00d714=1218\", \"51\u00d733=1683\")\n    ,@(\"79\u00d759=4661\", \"73\u00d726=1898\")\n    ,@(\"46\u00d762=2852\", \"77\u00d796=7392\")\n    ,@(\"29\u00d743=1247\", \"95\u00d759=5605\")\n    ,@(\"16\u00d775=1200\", \"77\u00d730=2310\")\n    ,@(\"24\u00d722=528\", \"52\u00d714=728\")\n    ,@(\"19\u00d771=1349\", \"65\u00d780=5200\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}"}
